$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '20.592.09'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.20%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.480.03'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.72%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.012'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.39%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9690'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.45%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '280.24'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.38%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3664'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.43%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3085'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.49%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '40.04'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.18%  '

$ws.Range("E10").Value = '  -0.17%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06666'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.37%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.006'
$ws.Range("D12").Style = "Normal"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.544'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.18%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '18.13'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.28%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.224'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.34%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9707'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.93%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001032'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.53%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.482.31'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.41%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.05967'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.90%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '70.03'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.27%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.509'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.40%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '14.53'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.96%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.08'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.65%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.270'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.60%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '20.629.33'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.41%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '142.68'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.40%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.118'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -8.37%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.30'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.77%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.645.21'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.34%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '114.15'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.26%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.981'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.60%  '

$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.8243'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.43%  '

$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.034'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.58%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08010'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.14%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.529'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.41%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.226'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +9.04%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.05816'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.97%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.757'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.88%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.9701'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.78%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.02049'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.24%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.638'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.21%  '

$ws.Range("E42").Value = '  -3.04%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1883'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.25%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5319'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.89%  '

$ws.Range("B45").Value = 'PancakeSwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.542'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.29%  '

$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.30'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.02%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '118.80'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.97%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5220'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.33%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.829'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.32%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06501'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.65%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.9902'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.36%  '
